$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.809.19"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "2.114.69"
$ws.Range("E3").Value = "  +6.53%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'333.28"
$ws.Range("E5").Value = "  +3.08%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +3.98%  "
$ws.Range("D8").Value = "'0.4380"
$ws.Range("E8").Value = "  +6.53%  "
$ws.Range("D9").Value = "'0.09003"
$ws.Range("E9").Value = "  +6.77%  "
$ws.Range("D10").Value = "'46.08"
$ws.Range("E10").Value = "  +8.29%  "
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("D12").Value = "'25.08"
$ws.Range("E12").Value = "  +3.80%  "
$ws.Range("D13").Value = "2.107.10"
$ws.Range("E13").Value = "  +7.19%  "
$ws.Range("D14").Value = "'6.758"
$ws.Range("E14").Value = "  +4.33%  "
$ws.Range("D15").Value = "'7.816"
$ws.Range("E15").Value = "  +5.67%  "
$ws.Range("D16").Value = "'97.43"
$ws.Range("E16").Value = "  +3.96%  "
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'0.00001127"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").Value = "'0.06662"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'6.345"
$ws.Range("E22").Value = "  +4.10%  "
$ws.Range("D23").Value = "30.857.27"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").Value = "'12.36"
$ws.Range("E24").Value = "  +7.52%  "
$ws.Range("D25").Value = "2.360.55"
$ws.Range("E25").Value = "  +7.50%  "
$ws.Range("D26").Value = "'2.271"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("D27").Value = "'22.78"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("E28").Value = "  +8.05%  "
$ws.Range("D29").Value = "'162.79"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'133.34"
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").Value = "'1.171"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").Value = "'6.224"
$ws.Range("E33").Value = "  +3.26%  "
$ws.Range("D34").Value = "'4.017"
$ws.Range("E34").Value = "  +5.38%  "
$ws.Range("E35").Value = "  +17.05%  "
$ws.Range("D36").Value = "'0.02614"
$ws.Range("E36").Value = "  +5.68%  "
$ws.Range("D37").Value = "'12.88"
$ws.Range("E37").Value = "  +9.13%  "
$ws.Range("D38").Value = "'5.533"
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("D39").Value = "'9.538"
$ws.Range("E39").Value = "  +6.76%  "
$ws.Range("D40").Value = "'0.06734"
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("D41").Value = "'0.2285"
$ws.Range("E41").Value = "  +5.14%  "
$ws.Range("D42").Value = "'0.6854"
$ws.Range("E42").Value = "  +3.94%  "
$ws.Range("D43").Value = "'1.250"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "'0.6436"
$ws.Range("E44").Value = "  +5.18%  "
$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "'14.07"
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("D47").Value = "'2.230"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").Value = "'3.660"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +4.37%  "
$ws.Range("D50").Value = "'82.60"
$ws.Range("E50").Value = "  +3.62%  "
$ws.Range("D51").Value = "'120.68"
$ws.Range("E51").Value = "  -2.58%  "
